$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update simple property values on the Metadata sheet ---
$ws1.Range("B3").Value = "0.1.7"
$ws1.Range("B6").Value = "draft"
$ws1.Range("B8").Value = "2024-08-27T12:23:18-05:00"
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Insert a new "Jurisdiction" row after the second Contact row (row 11), ---
# --- pushing the Description/Purpose/Copyright/Immutable rows down by one.  ---
# We avoid Rows.Insert() (it allocates a brand-new, unused cell style) and
# instead extend formatting to the new last row, then shift the values down
# manually (bottom-up so nothing is overwritten before it is copied).

# Extend the existing "data row" formatting (currently used by rows 2-15) down
# to row 16 by copying the format of row 15 onto it.
$ws1.Range("A15:B15").Copy()
$ws1.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Shift the values for rows 12-15 down to rows 13-16 (bottom-up).
$ws1.Range("A16").Value = $ws1.Range("A15").Value()
$ws1.Range("B16").Value = $ws1.Range("B15").Value()

$ws1.Range("A15").Value = $ws1.Range("A14").Value()
$ws1.Range("B15").Value = $ws1.Range("B14").Value()

$ws1.Range("A14").Value = $ws1.Range("A13").Value()
$ws1.Range("B14").Value = $ws1.Range("B13").Value()

$ws1.Range("A13").Value = $ws1.Range("A12").Value()
$ws1.Range("B13").Value = $ws1.Range("B12").Value()

# Fill in the new "Jurisdiction" row.
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""

Write-Output "done"
